$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "235.43") are stored as text, matching the source data which
# keeps all Price/Volume columns as text (inlineStr) cells.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

# --- Cell value updates ---
$ws.Range("D2").Value = '30.220.37'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '1.854.03'
$ws.Range("E3").Value = '  -2.25%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '235.43'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4780'
$ws.Range("E7").Value = '  -2.80%  '
$ws.Range("E8").Value = '  -4.00%  '
$ws.Range("D9").Value = '0.06461'
$ws.Range("E9").Value = '  -3.43%  '
$ws.Range("D10").Value = '1.852.05'
$ws.Range("E10").Value = '  -2.42%  '
$ws.Range("D11").Value = '0.07373'
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("E12").Value = '  -4.10%  '
$ws.Range("D13").Value = '5.090'
$ws.Range("E13").Value = '  -1.61%  '
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").Value = '0.6450'
$ws.Range("E15").Value = '  -3.05%  '
$ws.Range("D16").Value = '30.160.95'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '13.13'
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("D19").Value = '0.000007554'
$ws.Range("E19").Value = '  -3.64%  '
$ws.Range("D20").Value = '227.23'
$ws.Range("E20").Value = '  +18.10%  '
$ws.Range("D21").Value = '2.096.28'
$ws.Range("E21").Value = '  -1.80%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '5.281'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("D24").Value = '6.081'
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("D25").Value = '9.195'
$ws.Range("E25").Value = '  -3.02%  '
$ws.Range("D26").Value = '163.59'
$ws.Range("E26").Value = '  +0.88%  '
$ws.Range("D27").Value = '18.50'
$ws.Range("E27").Value = '  +1.54%  '
$ws.Range("D28").Value = '1.920'
$ws.Range("E28").Value = '  -0.70%  '
$ws.Range("D29").Value = '1.440'
$ws.Range("E29").Value = '  -2.17%  '
$ws.Range("D30").Value = '0.09175'
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("D31").Value = '4.229'
$ws.Range("E31").Value = '  -2.30%  '
$ws.Range("D32").Value = '3.951'
$ws.Range("E32").Value = '  -2.42%  '
$ws.Range("D33").Value = '0.04969'
$ws.Range("E33").Value = '  -3.86%  '
$ws.Range("D34").Value = '0.7285'
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").Value = '1.140'
$ws.Range("E35").Value = '  +3.53%  '
$ws.Range("D36").Value = '2.690'
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("D37").Value = '0.01841'
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("D38").Value = '2.596'
$ws.Range("E38").Value = '  -3.07%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '0.8984'
$ws.Range("E39").Value = '  -2.75%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '2.040'
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").Value = '5.935'
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").Value = '105.89'
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("D43").Value = '1.0000'
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("D44").Value = '0.4228'
$ws.Range("E44").Value = '  -3.57%  '
$ws.Range("D45").Value = '7.355'
$ws.Range("E45").Value = '  -3.05%  '
$ws.Range("E46").Value = '  -4.04%  '
$ws.Range("D47").Value = '63.95'
$ws.Range("E47").Value = '  -6.83%  '
$ws.Range("D48").Value = '1.496'
$ws.Range("E48").Value = '  +6.89%  '
$ws.Range("D49").Value = '8.722'
$ws.Range("E49").Value = '  -3.35%  '
$ws.Range("D50").Value = '33.74'
$ws.Range("E50").Value = '  -3.30%  '
$ws.Range("D51").Value = '0.05653'
$ws.Range("E51").Value = '  -3.44%  '

# Restore the default ("Normal") style on the price column so no stray
# number-format style lingers on cells that did not need one.
$priceCol.Style = "Normal"
